# Update cryptos list (price + volume(1h)) per the Wed Apr 17 13:08:30 UTC 2024 GitHub Actions run.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-Row {
    param(
        [int]$Row,
        [string]$D,
        [string]$E
    )
    if ($D -ne "") { $ws.Range("D$Row").Value = $D }
    if ($E -ne "") { $ws.Range("E$Row").Value = $E }
}

# Row 2 - Bitcoin
Set-Row 2 "62.421.06" "  -0.95%  "
# Row 3 - Ethereum
Set-Row 3 "3.037.19" "  -1.16%  "
# Row 4 - TetherUSD
Set-Row 4 "" "  -0.03%  "
# Row 5 - BNB
Set-Row 5 "537.29" "  -0.32%  "
# Row 6 - Solana
Set-Row 6 "132.72" "  -0.32%  "
# Row 7 - USDC
Set-Row 7 "0.999" "  +0.02%  "
# Row 8 - LidoStakedEther
Set-Row 8 "3.031.51" "  -1.16%  "
# Row 9 - XRP
Set-Row 9 "0.489" "  +0.24%  "
# Row 10 - Dogecoin
Set-Row 10 "0.153" "  -1.18%  "
# Row 11 - Toncoin
Set-Row 11 "6.17" "  -0.53%  "
# Row 12 - Cardano
Set-Row 12 "0.446" "  -2.44%  "
# Row 13 - ShibaInu
Set-Row 13 "0.0000219" "  -2.16%  "
# Row 14 - Avalanche
Set-Row 14 "33.79" "  -1.54%  "
# Row 15 - WrappedliquidstakedEther2.0
Set-Row 15 "3.519.28" "  -0.40%  "
# Row 16 - TRON
Set-Row 16 "" "  +1.46%  "
# Row 17 - WrappedBTC
Set-Row 17 "62.454.86" "  -0.76%  "
# Row 18 - WrappedEther
Set-Row 18 "3.030.61" "  -1.39%  "
# Row 19 - Polkadot
Set-Row 19 "6.54" "  -0.55%  "
# Row 20 - BitcoinCash
Set-Row 20 "463.84" "  -3.97%  "
# Row 21 - Chainlink
Set-Row 21 "13.25" "  +0.03%  "
# Row 22 - Polygon
Set-Row 22 "0.684" "  -2.27%  "
# Row 23 - Uniswap
Set-Row 23 "6.92" "  -3.21%  "
# Row 24 - Litecoin
Set-Row 24 "77.56" "  -0.91%  "
# Row 25 - InternetComputer(DFINITY)
Set-Row 25 "12.01" "  +0.09%  "
# Row 26 - Dai
Set-Row 26 "" "  -0.01%  "
# Row 27 - PancakeSwap
Set-Row 27 "2.67" "  -0.60%  "
# Row 28 - RenderToken
Set-Row 28 "7.74" "  -5.00%  "
# Row 29 - FirstDigitalUSD
Set-Row 29 "" "  -0.22%  "
# Row 30 - EthereumClassic
Set-Row 30 "25.77" "  -0.46%  "

# Row 31 / 32 swap: Mantle <-> ImmutableX (plus value updates)
$ws.Range("B31").Value = "ImmutableX"
$ws.Range("C31").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D31").Value = "1.85"
$ws.Range("E31").Value = "  -2.08%  "

$ws.Range("B32").Value = "Mantle"
$ws.Range("C32").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D32").Value = "1.14"
$ws.Range("E32").Value = "  +3.16%  "

# Row 33 - OKB
Set-Row 33 "58.36" "  -1.07%  "
# Row 34 - Stacks
Set-Row 34 "2.27" "  -5.74%  "
# Row 35 - NEARProtocol
Set-Row 35 "5.37" "  +4.54%  "
# Row 36 - Filecoin
Set-Row 36 "5.87" "  -1.27%  "
# Row 37 - Bittensor
Set-Row 37 "463.28" "  +0.87%  "
# Row 38 - Maker
Set-Row 38 "3.190.81" "  +2.31%  "
# Row 39 - VeChain
Set-Row 39 "0.0388" "  -0.21%  "
# Row 40 - Hedera
Set-Row 40 "0.0785" "  -0.16%  "
# Row 41 - Kaspa
Set-Row 41 "0.117" "  +2.20%  "
# Row 42 - Cosmos
Set-Row 42 "8.01" "  -0.02%  "
# Row 43 - dogwifhat (force text so the trailing zero in "2.50" survives)
$ws.Range("D43").NumberFormat = "@"
Set-Row 43 "2.50" "  -0.54%  "
# Row 44 - USDe
Set-Row 44 "" "  +0.15%  "
# Row 45 - TheGraph
Set-Row 45 "0.246" "  -0.78%  "

# Row 46 / 47 swap: Monero <-> InjectiveProtocol (plus value updates)
$ws.Range("B46").Value = "InjectiveProtocol"
$ws.Range("C46").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D46").Value = "24.94"
$ws.Range("E46").Value = "  +2.85%  "

$ws.Range("B47").Value = "Monero"
$ws.Range("C47").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D47").Value = "121.31"
$ws.Range("E47").Value = "  +3.07%  "

# Row 48 - Stellar
Set-Row 48 "0.108" "  +0.84%  "
# Row 49 - Fetch.AI
Set-Row 49 "" "  -2.28%  "
# Row 50 - PEPE
Set-Row 50 "0.0₃0510" "  +0.74%  "
# Row 51 - BitgetToken
Set-Row 51 "" "  +4.86%  "
